$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, row 1, matching the formatting
# of the other header cells (copy format from G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Make sure the value survived the format paste.
$ws.Range("H1").Value = "Save"

# Add the corresponding data value in row 2.
$ws.Range("H2").Value = 1
